$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (row 1): dataset labels for each propeller_demanded_0 column ---
$ws.Range("E1").Value = "Aframax"
$ws.Range("F1").Value = "Suezmax"
$ws.Range("G1").Value = "Aframax Casos Estranhos"
$ws.Range("H1").Value = "L280B50T17"

# --- Row 2 column headers for the new data columns ---
$ws.Range("F2").Value = "propeller_demanded_0"
$ws.Range("G2").Value = "propeller_demanded_0"

# --- Column F: Suezmax propeller_demanded_0 values ---
$ws.Range("F3").Value = -1420.0200000000002
$ws.Range("F4").Value = -858.48
$ws.Range("F5").Value = -631.12000000000012
$ws.Range("F6").Value = -438.06000000000006
$ws.Range("F7").Value = 0
$ws.Range("F8").Value = 730.1
$ws.Range("F9").Value = 1051.54
$ws.Range("F10").Value = 1431.78
$ws.Range("F11").Value = 2366.7000000000003
$ws.Range("F12").Value = 2921.3800000000006

# --- Column G: Aframax Casos Estranhos values (only a short, sparse range) ---
$ws.Range("G7").Value = 95
$ws.Range("G8").Value = 200
$ws.Range("G9").Value = 400
$ws.Range("G10").Value = 800

# --- Column H: L280B50T17 values ---
$ws.Range("H3").Value = -74
$ws.Range("H4").Value = -58
$ws.Range("H5").Value = -41
$ws.Range("H6").Value = -29
$ws.Range("H7").Value = 0
$ws.Range("H8").Value = 29
$ws.Range("H9").Value = 41
$ws.Range("H10").Value = 58
$ws.Range("H11").Value = 74
$ws.Range("H12").Value = 82

# --- Column widths: F & G match E's "bestFit" width; H is narrower ---
$ws.Columns(6).ColumnWidth = 21.5
$ws.Columns(7).ColumnWidth = 21.5
$ws.Columns(8).ColumnWidth = 10.17

# --- Final selection / active cell ---
[void]$ws.Range("H13").Select()

# --- Window geometry (best-effort; mirrors the author's resized Excel window) ---
$excel.ActiveWindow.Width = 24465
$excel.ActiveWindow.Height = 11445
